$wb = $excel.ActiveWorkbook
$evt = $wb.Worksheets.Item("Events")
$ws = $wb.Worksheets.Item("Incident")

# Headers (set values first so the format-only paste below doesn't clobber them)
$ws.Range("A1").Value = "IncName"
$ws.Range("B1").Value = "contractName"

# Re-use the same bold/yellow-fill header style already used on the "Events"
# sheet (cellXfs index 1) instead of building a brand-new style.
$evt.Range("A1:B1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

# Data rows
$ws.Range("A2").Value = "LDAP is Down.This is a Test Incident"
$ws.Range("B2").Value = "CAH_TEST"
$ws.Range("A3").Value = "Sentinal is Down.This is a Test Incident"
$ws.Range("B3").Value = "ES_TEST"

# Column widths (best-fit sized to content)
$ws.Columns.Item(1).ColumnWidth = 32.3333333
$ws.Columns.Item(2).ColumnWidth = 11.6666667

# Page setup to match the other sheet
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection on the sheet
$ws.Range("A3").Select()

# Make Incident the active sheet/tab
$ws.Activate()
